# Fix inconsistent naming for store fields in route normalization.
# This adjusts the generated SHIP_DISTANCES / SHIP_ROUTES sheets to add a
# new "Osborne"-related ship distance pair and a brand-new Route 9 (an FA
# route out of Gladstone landing at Melbourne), and bumps the Log sheet's
# "ensure_ship_distances_sheet" row-added counter to reflect the two new
# distance rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) SHIP_DISTANCES: insert two new "Osborne" distance rows.
#    Before:
#       9  Gladstone   Newcastle    674
#       10 Gladstone   Port Kembla  800
#       11 Gladstone   Townsville   433
#       12 Import_CL   Port Kembla  1000
#       13 Import_GBFS Port Kembla  1000
#       14 Melbourne   Port Kembla  637
#    After:
#       9  Gladstone   Newcastle    674
#       10 Gladstone   Osborne      (blank)
#       11 Gladstone   Port Kembla  800
#       12 Gladstone   Townsville   433
#       13 Import_CL   Port Kembla  1000
#       14 Import_GBFS Port Kembla  1000
#       15 Melbourne   Osborne      (blank)
#       16 Melbourne   Port Kembla  637
# ---------------------------------------------------------------------
$shipDist = $wb.Worksheets.Item("SHIP_DISTANCES")

# Insert from the bottom up so earlier row numbers stay stable.
$shipDist.Range("A14:C14").EntireRow.Insert()
$shipDist.Range("A14").Value = "Melbourne"
$shipDist.Range("B14").Value = "Osborne"

$shipDist.Range("A10:C10").EntireRow.Insert()
$shipDist.Range("A10").Value = "Gladstone"
$shipDist.Range("B10").Value = "Osborne"

# ---------------------------------------------------------------------
# 2) SHIP_ROUTES: insert a new "Route 9" column (between old Route 8 at
#    column I and old Route 9 at column J), an FA route Gladstone -> Melbourne,
#    and append a "Route 16" header at the end (column Q) for the old
#    last route's new slot.
# ---------------------------------------------------------------------
$shipRoutes = $wb.Worksheets.Item("SHIP_ROUTES")

# Shift every old Route 9..15 (columns J..P) one column to the right
# (K..Q), opening up a blank column J for the new route.
$shipRoutes.Range("J1").EntireColumn.Insert()

# Row 1 (route header labels) does not follow the generic shift: restore
# the simple "Route 1".."Route 16" sequence (old J1.."Route 9" must stay
# put, and "Route 16" is simply appended at the end in column Q).
$shipRoutes.Range("B1").Value = "Route 1"
$shipRoutes.Range("C1").Value = "Route 2"
$shipRoutes.Range("D1").Value = "Route 3"
$shipRoutes.Range("E1").Value = "Route 4"
$shipRoutes.Range("F1").Value = "Route 5"
$shipRoutes.Range("G1").Value = "Route 6"
$shipRoutes.Range("H1").Value = "Route 7"
$shipRoutes.Range("I1").Value = "Route 8"
$shipRoutes.Range("J1").Value = "Route 9"
$shipRoutes.Range("K1").Value = "Route 10"
$shipRoutes.Range("L1").Value = "Route 11"
$shipRoutes.Range("M1").Value = "Route 12"
$shipRoutes.Range("N1").Value = "Route 13"
$shipRoutes.Range("O1").Value = "Route 14"
$shipRoutes.Range("P1").Value = "Route 15"
$shipRoutes.Range("Q1").Value = "Route 16"

# New Route 9 (column J) field values - an FA route, Gladstone -> Melbourne.
$shipRoutes.Range("J2").Value = "North"
$shipRoutes.Range("J3").Value = 1.9
$shipRoutes.Range("J4").Value = "Gladstone"
$shipRoutes.Range("J5").Value = "FA"
$shipRoutes.Range("J6").Value = "FA_EXPSILO_STORE"
$shipRoutes.Range("J9").Value = "Melbourne"
$shipRoutes.Range("J10").Value = "FA"
$shipRoutes.Range("J11").Value = "FA_STORE"
$shipRoutes.Range("J21").Value = "Gladstone"

# Rows 14, 15 and 17 only ever had data in columns E (Route 4) and I
# (Route 8) - the generic column-insert above left their Route 8 value
# sitting in I; move it over to J (now Route 8's column stays I, with
# nothing new for Route 9) and clear the old spot.
# (Use .Value2 to read back a plain scalar instead of a COM Variant
# wrapper object.)
$shipRoutes.Range("J14").Value = $shipRoutes.Range("I14").Value2
$shipRoutes.Range("I14").Value = ""

$shipRoutes.Range("J15").Value = $shipRoutes.Range("I15").Value2
$shipRoutes.Range("I15").Value = ""

$shipRoutes.Range("J17").Value = $shipRoutes.Range("I17").Value2
$shipRoutes.Range("I17").Value = ""

# Rows 18, 19 and 20 gain a genuinely new Route 8 value (Osborne / FA /
# FA_STORE); the old Route 8 value shifts over to J.
$shipRoutes.Range("J18").Value = $shipRoutes.Range("I18").Value2
$shipRoutes.Range("I18").Value = "Osborne"

$shipRoutes.Range("J19").Value = $shipRoutes.Range("I19").Value2
$shipRoutes.Range("I19").Value = "FA"

$shipRoutes.Range("J20").Value = $shipRoutes.Range("I20").Value2
$shipRoutes.Range("I20").Value = "FA_STORE"

# ---------------------------------------------------------------------
# 3) Log: the "ensure_ship_distances_sheet" step now reports the two
#    rows added above.
# ---------------------------------------------------------------------
$log = $wb.Worksheets.Item("Log")
$log.Range("C33").Value = "rows_added=2"
